$wb = $excel.ActiveWorkbook

# "To Do- FY15 Release" gets a new task row inserted at row 16, and becomes the
# active/selected sheet (it was "SAM Variable Changes" before).
$ws = $wb.Worksheets.Item("To Do- FY15 Release")

$ws.Rows.Item(16).Insert()

$ws.Range("A16").Value = "Done"
$ws.Range("B16").Value = "Fix bug in performance ratio"
$ws.Range("C16").Value = "Janine"
$ws.Range("E16").Value = "A"

$ws.Activate()
$ws.Range("A17").Select()
